$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.973.59"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "1.815.33"
$ws.Range("E3").Value = "  +2.56%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.49%  "

$ws.Range("D5").Value = "'313.18"
$ws.Range("E5").Value = "  +2.24%  "

$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").Value = "'0.4298"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'0.3673"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").Value = "'0.07236"
$ws.Range("E9").Value = "  +0.56%  "

$ws.Range("D10").Value = "2.201.07"

$ws.Range("D11").Value = "'0.8635"
$ws.Range("E11").Value = "  +1.59%  "

$ws.Range("D12").Value = "'21.19"
$ws.Range("E12").Value = "  +4.07%  "

$ws.Range("D13").Value = "'5.400"
$ws.Range("E13").Value = "  +3.18%  "

$ws.Range("D14").Value = "'6.598"
$ws.Range("E14").Value = "  +2.64%  "

$ws.Range("D15").Value = "'0.06952"
$ws.Range("E15").Value = "  +0.45%  "

$ws.Range("D16").Value = "'81.16"
$ws.Range("E16").Value = "  +2.32%  "

$ws.Range("D17").Value = "'1.013"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").Value = "'0.000008909"
$ws.Range("E18").Value = "  +3.14%  "

$ws.Range("D19").Value = "'1.006"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").Value = "'15.19"
$ws.Range("E20").Value = "  +1.03%  "

$ws.Range("D21").Value = "27.027.73"
$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("D22").Value = "'5.168"
$ws.Range("E22").Value = "  +1.36%  "

$ws.Range("D23").Value = "2.451.37"
$ws.Range("E23").Value = "  +23.54%  "

$ws.Range("E24").Value = "  -1.82%  "

$ws.Range("D25").Value = "'153.86"

$ws.Range("D26").Value = "'1.879"
$ws.Range("E26").Value = "  +0.24%  "

$ws.Range("D27").Value = "'18.28"
$ws.Range("E27").Value = "  +1.23%  "

$ws.Range("D28").Value = "'5.220"
$ws.Range("E28").Value = "  +2.80%  "

$ws.Range("D29").Value = "'1.904"
$ws.Range("E29").Value = "  +9.01%  "

$ws.Range("D30").Value = "'114.54"

$ws.Range("D31").Value = "'0.08952"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("D32").Value = "'1.189"
$ws.Range("E32").Value = "  +6.96%  "

$ws.Range("D33").Value = "'0.7497"
$ws.Range("E33").Value = "  +3.55%  "

$ws.Range("D34").Value = "'4.421"
$ws.Range("E34").Value = "  +2.04%  "

$ws.Range("D35").Value = "'2.809"
$ws.Range("E35").Value = "  +2.34%  "

$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("D37").Value = "'1.133"
$ws.Range("E37").Value = "  +4.82%  "

$ws.Range("D38").Value = "'0.05210"
$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("D39").Value = "'0.01924"
$ws.Range("E39").Value = "  +1.93%  "

$ws.Range("D40").Value = "'0.5106"
$ws.Range("E40").Value = "  +3.78%  "

$ws.Range("D41").Value = "'2.750"
$ws.Range("E41").Value = "  +6.71%  "

$ws.Range("D42").Value = "'0.1647"
$ws.Range("E42").Value = "  +2.45%  "

$ws.Range("D43").Value = "'6.495"
$ws.Range("E43").Value = "  +3.66%  "

$ws.Range("D44").Value = "'8.322"
$ws.Range("E44").Value = "  +3.95%  "

$ws.Range("D45").Value = "'106.84"
$ws.Range("E45").Value = "  +1.93%  "

$ws.Range("D46").Value = "'10.36"
$ws.Range("E46").Value = "  +1.04%  "

$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").Value = "'0.4562"
$ws.Range("E48").Value = "  +1.56%  "

$ws.Range("D49").Value = "'1.645"
$ws.Range("E49").Value = "  +3.50%  "

$ws.Range("D50").Value = "'0.06212"
$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").Value = "'1.848"
$ws.Range("E51").Value = "  +6.20%  "
